$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

# Row 2
$ws.Range("C2").Value = 10.4

# Row 3
$ws.Range("B3").Value = 9.6
$ws.Range("D3").Value = 10.4
$ws.Range("F3").Value = 10.26

# Row 4
$ws.Range("C4").Value = 9.6
$ws.Range("E4").Value = 10.65
$ws.Range("F4").Value = 10.14

# Row 5
$ws.Range("D5").Value = 9.35
$ws.Range("F5").Value = 10.31
$ws.Range("G5").Value = 9.65
$ws.Range("I5").Value = 7

# Row 6
$ws.Range("C6").Value = 9.74
$ws.Range("D6").Value = 9.859999999999999
$ws.Range("E6").Value = 9.69
$ws.Range("G6").Value = 10.28
$ws.Range("H6").Value = 10.69

# Row 7
$ws.Range("E7").Value = 10.35
$ws.Range("F7").Value = 9.720000000000001
$ws.Range("H7").Value = 9.65

# Row 8
$ws.Range("F8").Value = 9.31
$ws.Range("G8").Value = 10.35
$ws.Range("I8").Value = 8.82

# Row 9
$ws.Range("E9").Value = 13
$ws.Range("H9").Value = 11.18
